$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Exp 14" experiment row (row 10)
$ws.Range("A10").Value = "Exp 14"
$ws.Range("B10").Value = 0.7
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = "Exp 14.png"

# Update the active selection to match the authored workbook state
$ws.Range("F14:F15").Select()
